$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row at 232 (pushes SUMBA GARCIA MARCOS ANTONIO and
# everything below it down by one row).
$ws1.Rows.Item(232).Insert()

# Fill in the new row's contents: same group (column A) as the
# surrounding rows, new client name in column B, zeros (matching the
# surrounding number format) in C:N.
$ws1.Range("A232").Value = "LOZANO MOLINA TITO"
$ws1.Range("B232").Value = "SOLORZANO BRAVO TERESA CONCEPCION"
for ($c = 3; $c -le 14; $c++) {
    $ws1.Cells.Item(232, $c).Value = 0
}

# The summary/totals row (previously row 260) is now row 261. Its
# "<n> de 258" labels need to become "<n> de 259" since the total
# number of clients increased by one.
for ($c = 3; $c -le 14; $c++) {
    $cell = $ws1.Cells.Item(261, $c)
    $cell.Value = $cell.Value2.Replace("de 258", "de 259")
}

# --- Sheet 2: "VENTA MENSUAL" -------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same insertion as sheet 1 (this sheet has no running "de NNN" labels,
# just numeric monthly totals, so nothing else needs patching after the
# insert shifts everything down).
$ws2.Rows.Item(232).Insert()

$ws2.Range("A232").Value = "LOZANO MOLINA TITO"
$ws2.Range("B232").Value = "SOLORZANO BRAVO TERESA CONCEPCION"
for ($c = 3; $c -le 6; $c++) {
    $ws2.Cells.Item(232, $c).Value = 0
}
